$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 800
$ws.Range("I18").Value = 800
$ws.Range("K18").Value = 800
$ws.Range("M18").Value = -516

$ws.Range("H41").Value = 765
$ws.Range("I41").Value = 687.6
$ws.Range("K41").Value = 687.6
$ws.Range("M41").Value = -247.6

$ws.Range("H58").Value = 10105.7
$ws.Range("I58").Value = 117.666664
$ws.Range("K58").Value = 352.999992
$ws.Range("M58").Value = -202.999992

$ws.Range("H98").Value = 1146.65
$ws.Range("I98").Value = 1181.2106
$ws.Range("K98").Value = 1181.2106
$ws.Range("M98").Value = 316.7893999999999

$ws.Range("H99").Value = 344
$ws.Range("I99").Value = 380
$ws.Range("J99").Value = 200
$ws.Range("K99").Value = 1140
$ws.Range("L99").Value = 600
$ws.Range("M99").Value = 358
$ws.Range("N99").Value = -3596

$ws.Range("H100").Value = 1810.1052
$ws.Range("I100").Value = 1990.375
$ws.Range("J100").Value = 848.6667
$ws.Range("K100").Value = 1990.375
$ws.Range("L100").Value = 848.6667
$ws.Range("M100").Value = -1449.375
$ws.Range("N100").Value = -1930.6667

$ws.Range("H122").Value = 1146.65
$ws.Range("I122").Value = 1181.2106
$ws.Range("K122").Value = 3543.6318
$ws.Range("M122").Value = -1093.6318

$ws.Range("H130").Value = 154483.75
$ws.Range("J130").Value = 154483.75
$ws.Range("L130").Value = 154483.75
$ws.Range("N130").Value = -164523.75

$ws.Range("H132").Value = 3013.889
$ws.Range("I132").Value = 2234.3572
$ws.Range("K132").Value = 6703.071599999999
$ws.Range("M132").Value = -4173.071599999999

$ws.Range("H139").Value = 162500
$ws.Range("I139").Value = 75000
$ws.Range("J139").Value = 250000
$ws.Range("K139").Value = 75000
$ws.Range("L139").Value = 250000
$ws.Range("M139").Value = -69860
$ws.Range("N139").Value = -260280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 12375.5
$ws.Range("I3").Value = 12375.5
$ws.Range("K3").Value = 12375.5
$ws.Range("M3").Value = -12260.5

$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 16
$ws.Range("N4").ClearContents()

$ws.Range("H5").Value = 65.875
$ws.Range("I5").Value = 47.833332
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 47.833332
$ws.Range("L5").Value = 120
$ws.Range("M5").Value = 64.166668
$ws.Range("N5").Value = -344

$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H12").Value = 199.33333
$ws.Range("I12").Value = 135.63637
$ws.Range("J12").Value = 900
$ws.Range("K12").Value = 135.63637
$ws.Range("L12").Value = 900
$ws.Range("M12").Value = 37.36363
$ws.Range("N12").Value = -1246

$ws.Range("H13").Value = 1696.5
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H16").Value = 3462.25
$ws.Range("I16").Value = 3462.25
$ws.Range("K16").Value = 3462.25
$ws.Range("M16").Value = -3175.25

$ws.Range("H22").Value = 6129.8
$ws.Range("I22").Value = 3633.3333
$ws.Range("K22").Value = 3633.3333
$ws.Range("M22").Value = -3334.3333

$ws.Range("H45").Value = 15791.5
$ws.Range("I45").Value = 15791.5
$ws.Range("K45").Value = 15791.5
$ws.Range("M45").Value = -15414.5

$ws.Range("H63").Value = 4966.8335
$ws.Range("J63").Value = 2998
$ws.Range("L63").Value = 2998
$ws.Range("N63").Value = -4370

$ws.Range("H66").Value = 4966.8335
$ws.Range("J66").Value = 2998
$ws.Range("L66").Value = 14990
$ws.Range("N66").Value = -21854

$ws.Range("H97").Value = 3758.0625
$ws.Range("I97").Value = 2087.077
$ws.Range("J97").Value = 10999
$ws.Range("K97").Value = 2087.077
$ws.Range("L97").Value = 10999
$ws.Range("M97").Value = -1591.077
$ws.Range("N97").Value = -11991

$ws.Range("H139").Value = 89707.5
$ws.Range("J139").Value = 89707.5
$ws.Range("L139").Value = 89707.5
$ws.Range("N139").Value = -99987.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 65.875
$ws.Range("I4").Value = 47.833332
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 47.833332
$ws.Range("L4").Value = 120
$ws.Range("M4").Value = 67.166668
$ws.Range("N4").Value = -350

$ws.Range("H22").Value = 20399.5
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 40299
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 40299
$ws.Range("M22").Value = -327
$ws.Range("N22").Value = -40645

$ws.Range("H25").Value = 18671.666
$ws.Range("J25").Value = 19008
$ws.Range("L25").Value = 19008
$ws.Range("N25").Value = -19478

$ws.Range("H29").Value = 7133.3335
$ws.Range("I29").Value = 10197.286
$ws.Range("J29").Value = 2843.8
$ws.Range("K29").Value = 10197.286
$ws.Range("L29").Value = 2843.8
$ws.Range("M29").Value = -9908.286
$ws.Range("N29").Value = -3421.8

$ws.Range("H94").Value = 2330.4614
$ws.Range("J94").Value = 3816.3333
$ws.Range("L94").Value = 3816.3333
$ws.Range("N94").Value = -4718.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2459.9333
$ws.Range("I105").Value = 1380
$ws.Range("J105").Value = 2999.9
$ws.Range("K105").Value = 1380
$ws.Range("L105").Value = 2999.9
$ws.Range("M105").Value = 367
$ws.Range("N105").Value = -6493.9

$ws.Range("H132").Value = 3060.3684
$ws.Range("I132").Value = 3060.3684
$ws.Range("K132").Value = 9181.1052
$ws.Range("M132").Value = -6651.1052

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 136.66667
$ws.Range("I61").Value = 105
$ws.Range("K61").Value = 315
$ws.Range("M61").Value = -100

$ws.Range("H92").Value = 538.9
$ws.Range("I92").Value = 202
$ws.Range("J92").Value = 576.3333
$ws.Range("K92").Value = 606
$ws.Range("L92").Value = 1728.9999
$ws.Range("M92").Value = 642
$ws.Range("N92").Value = -4224.9999

$ws.Range("H131").Value = 43354.957
$ws.Range("J131").Value = 1918.8823
$ws.Range("L131").Value = 5756.6469
$ws.Range("N131").Value = -15836.6469

$ws.Range("H137").Value = 2775.1667
$ws.Range("J137").Value = 3424.625
$ws.Range("L137").Value = 10273.875
$ws.Range("N137").Value = -20473.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 170.72223
$ws.Range("I2").Value = 121.94118
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 121.94118
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -8.941180000000003
$ws.Range("N2").Value = -1226

$ws.Range("H11").Value = 9000500
$ws.Range("I11").Value = 9000500
$ws.Range("K11").Value = 9000500
$ws.Range("M11").Value = -9000361

$ws.Range("H18").Value = 37072336
$ws.Range("I18").Value = 55558504
$ws.Range("K18").Value = 55558504
$ws.Range("M18").Value = -55558211

$ws.Range("H102").Value = 31254890
$ws.Range("I102").Value = 1878.4828
$ws.Range("K102").Value = 1878.4828
$ws.Range("M102").Value = -256.4828

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3579
$ws.Range("I68").Value = 3590.2856
$ws.Range("J68").Value = 3500
$ws.Range("K68").Value = 3590.2856
$ws.Range("L68").Value = 3500
$ws.Range("M68").Value = -2841.2856
$ws.Range("N68").Value = -4998

$ws.Range("H71").Value = 3579
$ws.Range("I71").Value = 3590.2856
$ws.Range("J71").Value = 3500
$ws.Range("K71").Value = 17951.428
$ws.Range("L71").Value = 17500
$ws.Range("M71").Value = -14207.428
$ws.Range("N71").Value = -24988

$ws.Range("H93").Value = 1539.4117
$ws.Range("I93").Value = 1544.6666
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 1544.6666
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = -296.6666
$ws.Range("N93").Value = -3996

$ws.Range("H104").Value = 21666.334
$ws.Range("J104").Value = 21666.334
$ws.Range("L104").Value = 21666.334
$ws.Range("N104").Value = -28654.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 7514833
$ws.Range("I15").Value = 9014000
$ws.Range("J15").Value = 18997
$ws.Range("K15").Value = 9014000
$ws.Range("L15").Value = 18997
$ws.Range("M15").Value = -9013712
$ws.Range("N15").Value = -19573

$ws.Range("H62").Value = 5899.2
$ws.Range("I62").Value = 4268.6
$ws.Range("J62").Value = 7529.8
$ws.Range("K62").Value = 4268.6
$ws.Range("L62").Value = 7529.8
$ws.Range("M62").Value = -3644.6
$ws.Range("N62").Value = -8777.799999999999

$ws.Range("H65").Value = 5899.2
$ws.Range("I65").Value = 4268.6
$ws.Range("J65").Value = 7529.8
$ws.Range("K65").Value = 21343
$ws.Range("L65").Value = 37649
$ws.Range("M65").Value = -18223
$ws.Range("N65").Value = -43889

$ws.Range("H122").Value = 2122.0417
$ws.Range("I122").Value = 1840.35
$ws.Range("K122").Value = 5521.049999999999
$ws.Range("M122").Value = -3071.049999999999

$ws.Range("H126").Value = 12689.667
$ws.Range("I126").Value = 15247.333
$ws.Range("J126").Value = 5016.6665
$ws.Range("K126").Value = 45741.999
$ws.Range("L126").Value = 15049.9995
$ws.Range("M126").Value = -43271.999
$ws.Range("N126").Value = -19989.9995

$ws.Range("H132").Value = 1882.6471
$ws.Range("I132").Value = 1882.6471
$ws.Range("K132").Value = 5647.9413
$ws.Range("M132").Value = -3117.9413
